# Update attendance/interest numbers ("想去人数") on the 展览 and 全部类型
# sheets to reflect the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 12922
$wsExhibition.Range("F7").Value = 45
$wsExhibition.Range("F10").Value = 12874
$wsExhibition.Range("F11").Value = 286
$wsExhibition.Range("F13").Value = 8694
$wsExhibition.Range("F14").Value = 7693
$wsExhibition.Range("F15").Value = 198
$wsExhibition.Range("F18").Value = 127
$wsExhibition.Range("F22").Value = 382
$wsExhibition.Range("F23").Value = 184
$wsExhibition.Range("F24").Value = 19

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 12922
$wsAllTypes.Range("F8").Value = 45
$wsAllTypes.Range("F11").Value = 12874
$wsAllTypes.Range("F12").Value = 286
$wsAllTypes.Range("F14").Value = 8694
$wsAllTypes.Range("F15").Value = 7693
$wsAllTypes.Range("F16").Value = 198
$wsAllTypes.Range("F19").Value = 127
$wsAllTypes.Range("F24").Value = 382
$wsAllTypes.Range("F25").Value = 184
$wsAllTypes.Range("F26").Value = 19
